$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 224 holds a date serial value
# that was updated by one day (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04).
$ws.Range("C2:C224").Value = 45203
